# "Allocation Web link Asset Web"
# Refresh the example/sample data row (row 6) on the "Ds cấp phát TSCĐ"
# sheet so the sample allocation date, asset code, and warehouse/source
# code reflect a current, valid example (linking to the live Asset/
# Allocation web records) instead of the stale placeholder values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ds cấp phát TSCĐ")

$ws.Range("B6").Value = "05/05/2025"
$ws.Range("C6").Value = "TS-008520"
$ws.Range("D6").Value = "XEXqJ1"

# Match the author's final on-screen selection (rows 8-17) when the
# workbook was saved.
$ws.Activate() | Out-Null
$ws.Rows("8:17").Select() | Out-Null
